$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the values from column B (rows 1-7) to column A (rows 1-7)
for ($r = 1; $r -le 7; $r++) {
    $val = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value2 = $val
    $ws.Cells.Item($r, 2).ClearContents()
}
